# Act worksheet: replace the generic "act one" placeholder inside the
# short-description paragraphs that actually belong to Act Two / Act
# Three with the correct act word, splitting the run the way Word does
# when the word is retyped in place (prefix run / replaced word run /
# suffix run, all keeping the original italic Calibri run formatting).

$d = $word.ActiveDocument

$placeholder = "A short description of what happens in act one.  An overview, if you will"

# Locate every paragraph whose text is this (still-generic) placeholder.
$targets = @()
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "A short description of what happens in act one.*") {
        $targets += $i
    }
}

foreach ($idx in $targets) {

    # Walk backwards to the nearest "Act <Word>:" heading that governs
    # this paragraph, so we know which word the placeholder should become.
    $newWord = $null
    for ($j = $idx - 1; $j -ge 1; $j--) {
        $headingText = $d.Paragraphs.Item($j).Range.Text
        if ($headingText -like "Act Two:*") { $newWord = "two"; break }
        if ($headingText -like "Act Three:*") { $newWord = "three"; break }
        if ($headingText -like "Act One:*") { $newWord = $null; break }
    }

    # Only rewrite the paragraphs that actually belong to Act Two / Act
    # Three - the one sitting right under "Act One:" keeps saying "one".
    if ($newWord -eq $null) {
        continue
    }

    $para = $d.Paragraphs.Item($idx)
    $pText = $para.Range.Text
    # Drop the trailing paragraph mark captured by Range.Text.
    $plain = $pText.TrimEnd([char]13, [char]7)

    $splitAt = $plain.IndexOf("one")
    $prefix = $plain.Substring(0, $splitAt)
    $suffix = $plain.Substring($splitAt + 3)

    $prefixXml = $prefix -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $newWordXml = $newWord -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $suffixXml = $suffix -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

    # Only mark xml:space="preserve" where the run text actually has
    # leading/trailing whitespace that needs protecting.
    $prefixTag = "<w:t>" + $prefixXml + "</w:t>"
    if ($prefix.Trim() -ne $prefix) {
        $prefixTag = '<w:t xml:space="preserve">' + $prefixXml + '</w:t>'
    }
    $newWordTag = "<w:t>" + $newWordXml + "</w:t>"
    if ($newWord.Trim() -ne $newWord) {
        $newWordTag = '<w:t xml:space="preserve">' + $newWordXml + '</w:t>'
    }
    $suffixTag = "<w:t>" + $suffixXml + "</w:t>"
    if ($suffix.Trim() -ne $suffix) {
        $suffixTag = '<w:t xml:space="preserve">' + $suffixXml + '</w:t>'
    }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body>' +
      '<w:p>' +
        '<w:pPr>' +
          '<w:pStyle w:val="Normal"/>' +
          '<w:bidi w:val="0"/>' +
          '<w:jc w:val="left"/>' +
          '<w:rPr/>' +
        '</w:pPr>' +
        '<w:r>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>' +
            '<w:i/>' +
            '<w:iCs/>' +
          '</w:rPr>' +
          $prefixTag +
        '</w:r>' +
        '<w:r>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>' +
            '<w:i/>' +
            '<w:iCs/>' +
          '</w:rPr>' +
          $newWordTag +
        '</w:r>' +
        '<w:r>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>' +
            '<w:i/>' +
            '<w:iCs/>' +
          '</w:rPr>' +
          $suffixTag +
        '</w:r>' +
      '</w:p>' +
      '</w:body>' +
      '</w:document>' +
      '</pkg:xmlData>' +
      '</pkg:part>' +
      '</pkg:package>'

    $para.Range.InsertXML($xml)
}
